$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (type_user) for row 2 first, so "Admin" becomes the first new
#     shared string introduced (matches the author's original entry order). ---
$ws.Cells.Item(2, 6).Value = "Admin"

# --- Column A (user name / "name" field) for the new rows, top to bottom. ---
$ws.Cells.Item(2, 1).Value = "ttestteo4"
$ws.Cells.Item(3, 1).Value = "ttestteo5"
$ws.Cells.Item(4, 1).Value = "ttestteo6"
$ws.Cells.Item(5, 1).Value = "ttestteo7"
$ws.Cells.Item(6, 1).Value = "ttestteo8"
$ws.Cells.Item(7, 1).Value = "ttestteo9"
$ws.Cells.Item(8, 1).Value = "ttestteo10"

# --- Column D (user_name) for the new rows, top to bottom. ---
$ws.Cells.Item(2, 4).Value = "testeusertest4"
$ws.Cells.Item(3, 4).Value = "testeusertest5"
$ws.Cells.Item(4, 4).Value = "testeusertest6"
$ws.Cells.Item(5, 4).Value = "testeusertest7"
$ws.Cells.Item(6, 4).Value = "testeusertest8"
$ws.Cells.Item(7, 4).Value = "testeusertest9"
$ws.Cells.Item(8, 4).Value = "testeusertest10"

# --- Remaining column F (type_user) values, alternating Admin/Estudiante. ---
$ws.Cells.Item(3, 6).Value = "Estudiante"
$ws.Cells.Item(4, 6).Value = "Admin"
$ws.Cells.Item(5, 6).Value = "Estudiante"
$ws.Cells.Item(6, 6).Value = "Admin"
$ws.Cells.Item(7, 6).Value = "Estudiante"
$ws.Cells.Item(8, 6).Value = "Admin"

# --- Fill the remaining repeated columns (B, C, E, G, H, I, J) for rows 3-8,
#     reusing the same values already present in row 2. ---
for ($r = 3; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value = "prueba"
    $ws.Cells.Item($r, 3).Value = 123456
    $ws.Cells.Item($r, 5).Value = "brayantriana22@gmail.com"
    $ws.Cells.Item($r, 7).Value = "calle siempreviva 123"
    $ws.Cells.Item($r, 8).Value = "https://www.google.com/url?sa=i&url=https%3A%2F%2Fes.123rf.com%2Fphoto_59070200_icono-de-usuario-hombre-perfil-hombre-de-negocios-avatar-icono-persona-en-la-ilustraci%25C3%25B3n-vectorial.html&psig=AOvVaw0KE_280JdOEhHeKpuGBrjB&ust=1645566286766000&source=images&cd=vfe&ved=0CAsQjRxqFwoTCIizhqPikfYCFQAAAAAdAAAAABAD"
    $ws.Cells.Item($r, 9).Value = 3007819686
    $ws.Cells.Item($r, 10).Value = 12345679987
}

# --- Hyperlinks for the new email cells, E3:E8 (E2's link already exists).
#     Restore the same cell style E2 already carries afterwards, since
#     Hyperlinks.Add otherwise stamps its own (slightly different) style. ---
for ($r = 3; $r -le 8; $r++) {
    $ws.Hyperlinks.Add($ws.Range("E$r"), "mailto:brayantriana22@gmail.com") | Out-Null
    $ws.Range("E$r").Style = $ws.Range("E2").Style
}

# --- Update the selection to match the saved view state. ---
$ws.Range("C10").Select()
